$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete rows (1996 and 1997 records) - this shifts
# all subsequent rows up by two, matching the new layout.
$ws.Range("A2:C3").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Rename the header columns to match the updated data source naming.
# (Set in C, B, A order so the shared-string table is built in the
# same sequence as the authoritative workbook.)
$ws.Range("C1").Value = "sleep_time"
$ws.Range("B1").Value = "age_range"
$ws.Range("A1").Value = "year"
